$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27
$ws.Range("N27").Value = 1.85
$ws.Range("O27").Value = 1.95

# Row 28
$ws.Range("G28").Value = 2.15
$ws.Range("H28").Value = 3.2
$ws.Range("I28").Value = 3
$ws.Range("J28").Value = 1.1
$ws.Range("K28").Value = 7
$ws.Range("L28").Value = 1.5
$ws.Range("M28").Value = 2.5
$ws.Range("N28").Value = 2.4
$ws.Range("O28").Value = 1.53
$ws.Range("P28").Value = 1.53
$ws.Range("Q28").Value = 2.38
$ws.Range("R28").Value = 2.1
$ws.Range("S28").Value = 1.67
$ws.Range("T28").Value = 6
$ws.Range("U28").Value = 9.5
$ws.Range("V28").Value = 10
$ws.Range("W28").Value = 21
$ws.Range("X28").Value = 21
$ws.Range("Y28").Value = 41
$ws.Range("Z28").Value = 7
$ws.Range("AA28").Value = 6.5
$ws.Range("AB28").Value = 21
$ws.Range("AC28").Value = 81
$ws.Range("AD28").Value = 101
$ws.Range("AE28").Value = 7.5
$ws.Range("AF28").Value = 15
$ws.Range("AG28").Value = 12
$ws.Range("AH28").Value = 41
$ws.Range("AI28").Value = 34
$ws.Range("AJ28").Value = 41

# Row 32
$ws.Range("G32").Value = 2.32
$ws.Range("H32").Value = 2.45
$ws.Range("I32").Value = 4.15
$ws.Range("K32").Value = 4.15
$ws.Range("N32").Value = 3.1
$ws.Range("P32").Value = 1.72
$ws.Range("S32").Value = 1.53
$ws.Range("U32").Value = 9.5
$ws.Range("W32").Value = 26
$ws.Range("Z32").Value = 4.15
$ws.Range("AE32").Value = 7.5

# Row 37
$ws.Range("G37").Value = 1.75
$ws.Range("H37").Value = 3.2
$ws.Range("I37").Value = 4.45
$ws.Range("L37").Value = 1.37
$ws.Range("M37").Value = 2.85
$ws.Range("N37").Value = 2.1
$ws.Range("O37").Value = 1.57
$ws.Range("P37").Value = 1.42
$ws.Range("Q37").Value = 2.35
$ws.Range("R37").Value = 2.03
$ws.Range("S37").Value = 1.7
$ws.Range("T37").Value = 4.9
$ws.Range("U37").Value = 6.3
$ws.Range("W37").Value = 11.25
$ws.Range("X37").Value = 12.5
$ws.Range("Z37").Value = 7.6
$ws.Range("AA37").Value = 5.6
$ws.Range("AB37").Value = 14.5
$ws.Range("AC37").Value = 70
$ws.Range("AE37").Value = 8.75
$ws.Range("AF37").Value = 19
$ws.Range("AG37").Value = 12.5
$ws.Range("AH37").Value = 60
$ws.Range("AI37").Value = 37
$ws.Range("AJ37").Value = 45

# Row 38
$ws.Range("G38").Value = 3.05
$ws.Range("H38").Value = 3
$ws.Range("I38").Value = 2.22
$ws.Range("N38").Value = 2.02
$ws.Range("O38").Value = 1.62
$ws.Range("P38").Value = 1.39
$ws.Range("Q38").Value = 2.42
$ws.Range("R38").Value = 1.86
$ws.Range("S38").Value = 1.84
$ws.Range("T38").Value = 7.5
$ws.Range("U38").Value = 13
$ws.Range("V38").Value = 9.25
$ws.Range("W38").Value = 32
$ws.Range("X38").Value = 22
$ws.Range("Y38").Value = 28
$ws.Range("Z38").Value = 8
$ws.Range("AA38").Value = 5.2
$ws.Range("AB38").Value = 11.5
$ws.Range("AC38").Value = 50
$ws.Range("AD38").Value = 350
$ws.Range("AE38").Value = 6
$ws.Range("AF38").Value = 8.75
$ws.Range("AG38").Value = 7.6
$ws.Range("AH38").Value = 17.5
$ws.Range("AI38").Value = 15.5
$ws.Range("AJ38").Value = 24

# Row 45
$ws.Range("I45").Value = 3
$ws.Range("L45").Value = 1.25
$ws.Range("N45").Value = 1.75
$ws.Range("O45").Value = 1.87
$ws.Range("AA45").Value = 6.6

# Row 57
$ws.Range("N57").Value = 1.75
$ws.Range("O57").Value = 2.05

# Row 71
$ws.Range("J71").Value = 1.08
$ws.Range("K71").Value = 8

# Row 81
$ws.Range("J81").Value = 1.03
$ws.Range("L81").Value = 1.22

# Row 83
$ws.Range("G83").Value = 2.1
$ws.Range("I83").Value = 3.3
$ws.Range("J83").Value = 1.04
$ws.Range("L83").Value = 1.33
$ws.Range("T83").Value = 7
$ws.Range("U83").Value = 9.5
$ws.Range("X83").Value = 19
